$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "U" value in column D for data rows 2-6, and clear the old
# F/G columns, keeping only a single "1" in column E for each row.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = $null
    $ws.Cells.Item($r, 7).Value = $null
}

# Update the active selection to D7, matching the saved selection state.
$ws.Range("D7").Select()
